# Merge the "Untitled Section" into the preceding "Default Section".
#
# Before:
#   Default Section  {BC24336F-...} -> slides 268, 277
#   Untitled Section {FFBA3020-...} -> slides 271, 278, 281, 257
#
# After:
#   Default Section  {BC24336F-...} -> slides 268, 277, 271, 278, 281, 257
#
# Deleting a section with DeleteSlides=$false removes the section marker
# but keeps its slides, which are absorbed into the previous section's
# slide list (matches PowerPoint's native "Remove Section" behavior).

$p = $ppt.ActivePresentation
$sections = $p.SectionProperties

for ($i = 1; $i -le $sections.Count; $i++) {
    if ($sections.Name($i) -eq "Untitled Section") {
        $sections.Delete($i, $false)
        break
    }
}
